$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shift the names up: A2->clinton, A3->bush, A4->rigan, and rename
# "coolige" to "Jimy" moving it to the bottom of the list (A5).
$ws.Range("A2").Value = "clinton"
$ws.Range("A3").Value = "bush"
$ws.Range("A4").Value = "rigan"
$ws.Range("A5").Value = "Jimy"

# A2 used to carry the Hyperlink cell style (for "coolige"); that
# name has moved, so restore the default "Normal" style on A2.
$ws.Range("A2").Style = "Normal"

# The Hyperlink cell style is no longer used anywhere in the workbook
# now that A2 is back to "Normal" - remove it so it isn't left behind
# as dead formatting.
$wb.Styles.Item("Hyperlink").Delete()

# Update the active selection to B15.
$ws.Range("B15").Select()
